$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 17 game-by-game values to long strings ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = '10 0 4 3 10 7 11 8 4 8 0 0 6 2 4 3 7 16 4 2 6 1 5 3 2 1 2 6 8 2 1 1 1 0 14 0 9 6 -1 6 2 8 12 8 9 3 4 6 3 1 6 6 -3 4 2 5 1 1 2 0 4 7 6 6 1 9 4 2 4 4 4 0 8 11 4 0 0 -1 3 1 1 3 -3 19 -1 0 10 -2 10 -1 18 8 1 3 4 12 9 3 3 5 0 4 1 1 1 5 4 4 2 1 0 3 54 7 4 14 0 11 3 -3 1 2 4 5 26 1 3 9 1 5 0 12 1 0 -2 11 4 1 2 10 1 2 8 4 4 40 3 5 0 9 4 0 1 0 -5 3 5 2 4 4 3 0 5 1 4 -2 3 3 -1 9 2 16 3 7 -2 3 7 0 14 -2 4 5 11 6 2 5 4 17 1 -1 6 1 4 -1 0 7 -2 6 16 9 2 1 18 2 7 12 0 1 4 -1 8 7 4 2 3 1 11 -1 4 5 4 2 2 2 7 6 9 0 0 9 1 3 3 8 1 4 1 4 2 1 1 2 3 3 1 -4 1 5 8 2 35 21 5 1 7 6 2 4 1 3 5 2 3 -3 2 6 -2 4 0 -1 1 6 2 4 0 4 0 0 4 4 6 11 11 1 9 0 7 0 8 4 2 1 2 2 5 2 0 7 -1 3 0 0 6 1 2 5 4 2 7 6 3 11 4 1 2 4 1 8 5 0 9 8 21 2 7 1 2 -4 5 10 7 2 14 -1 -6 10 2 5 4 9 0 -6 9 3 3 6 3 0 -4 2 13 1 1 0 0 1 7 2 5 12 8 16 0 17 6 2 1 8 1 1 0 6 2 5 1 22 0 3 23 6 1 2 1 6 3 0 4 12 6 5 4 2 2 10 57 0 2 1 4 1 0 5 1 7 0 11 10 2 24 14 1 2 1 2 2 1 -5 2 7 4 1 -3 2 7 12 5 0 1 2 2 1 1 15 3 11 6 9'
$ydsWs.Range("C2").Value = '-3 -1 2 7 -1 2 2 0 0 1 16 1 19 4 8 8 -1 4 7 8 1 2 11 2 16 14 6 7 3 0 15 2 0 29 9 6 12 0 10 0 13 2 3 0 -2 8 1 1 9 1 6 -1 5 4 -1 2 -2 4 -2 4 5 4 18 -1 6 17 8 10 3 21 6 6 4 4 4 10 9 5 6 3 7 15 2 9 15 9 2 0 3 3 6 4 11 3 8 5 7 0 6 0 14 0 2 1 -3 5 11 9 9 7 7 14 4 -2 1 7 7 12 1 13 -2 15 1 5 6 7 0 5 2 2 -2 4 17 11 0 5 -1 -3 9 3 2 11 3 3 2 4 1 4 4 -2 -2 2 4 13 3 1 6 2 2 0 2 3 7 4 12 6 5 8 2 1 11 4 1 3 8 6 0 10 -2 7 0 9 4 1 4 2 4 8 16 0 0 4 9 -5 -2 11 14 5 2 -1 2 3 1 -1 2 11 11 12 10 6 2 4 12 1 16 15 10 2 5 4 1 1 4 0 5 2 10 4 -9 3 15 4 8 6 -4 -3 1 7 -2 18 4 5 2 3 11 4 0 4 9 2 0 1 2 -4 7 15 2 1 0 4 9 2 8 7 3 1 13 0 24 -1 7 -2 1 2 2 3 3 8 1 0 5 5 4 6 4 3 19 -1 11 12 4 0 -1 2 5 2 2 0 0 10 3 6 6 0 2 10 0 7 6 -3 5 5 0 27 6 8 0 0 -3 1 0 1 4 5 3 0 6 0 2 8 1 16 5 0 2 1 2 -2 3 43 3 -1 4 1 -2 2 2 4 -2 7 -2 5 5 11 -1 2 8 18 4 3 -1 1 3 3 -1 3 2 2 1 -3 3 4 -1 2 0 -1 5 4 6 3 3 1 0 12 -1 4 1'
$ydsWs.Range("B3").Value = '38 1 21 8 5 3 10 3 17 18 4 22 3 12 36 3 0 10 26 11 9 4 8 15 16 4 0 15 3 34 9 19 77 4 18 6 13 25 9 11 29 6 2 6 17 2 3 7 7 7 5 28 9 7 8 6 10 10 -2 12 5 6 34 8 7 9 18 33 13 1 36 3 7 19 7 5 8 25 5 41 6 14 10 7 9 -3 14 7 10 3 16 9 7 9 18 6 12 25 5 7 8 13 10 9 23 33 3 6 2 9 14 8 5 8 13 7 9 2 30 9 11 0 15 21 10 8 5 13 -1 10 33 5 0 13 1 34 8 9 20 11 14 13 11 1 8 10 6 1 41 5 19 15 16 2 11 12 13 47 9 8 4 8 58 6 9 8 3 6 8 11 20 9 9 4 9 23 29 9 6 22 10 0 9 9 8 0 21 17 33 50 17 6 -2 25 45 4 4 1 5 5 9 6 0 1 4 19 5 20 4 0 5 19 21 7 9 -5 17 5 3 4 8 17 5 7 11 1 -5 31 7 1 12 25 2 10 20 6 11 4 21 6 -2 14 7 5 16 5 3 10 9 15 7 3 12 20 10 12 23 8 7 12 13 0 4 14 14 6 17 11 14 5 7 2 18 -2 26 10 41 13 3 9 0 5 12 6 6 7 23 15 36 17 12 8 -2 3 11 9 6 -1 14 13 2 5 24 10 15 0 3 9 20 19 27 16 5 7 6 7 13 11 26 3 10 5 10 5 8 13 7 11 4 3 11 10 24 3 33 4 5 7 16 7 4 7 14 2 14 8 42 4 17 3 3 3 20 10 5 3 4 9 4 5 5 11 43 -7 6 19 16 7 12 7 11'
$ydsWs.Range("C3").Value = '3 4 8 12 11 39 10 9 9 4 4 19 13 9 10 5 10 8 10 8 7 64 6 14 5 3 7 9 14 -3 4 15 6 6 12 13 7 13 3 17 12 9 8 6 0 8 14 11 20 7 9 12 13 5 12 11 6 12 2 12 9 13 6 9 20 5 18 5 2 7 14 14 7 35 18 9 0 6 8 7 22 8 6 7 13 9 12 5 9 14 16 14 26 6 1 18 6 12 4 26 5 13 26 13 4 25 20 4 17 10 5 11 5 7 6 57 24 13 3 6 6 7 6 3 10 10 0 22 5 5 0 -1 18 8 6 5 2 4 5 -3 1 2 10 8 10 5 9 7 7 12 4 1 4 15 -2 7 2 19 33 8 5 6 9 12 10 4 7 4 19 18 16 7 11 14 22 11 10 16 6 5 19 29 1 11 32 7 10 15 6 5 12 5 3 7 5 13 17 2 6 16 -2 6 10 2 3 11 6 9 2 10 6 10 1 2 3 9 8 15 7 3 6 9 25 14 36 4 12 16 6 16 48 2 10 3 12 3 5 6 3 10 7 6 10 34 11 2 1 7 12 11 6 6 2 7 19 0 11 5 13 1 7 11 8 2 17 6 40 6 2 6 52 4 4 23 4 17 44 7 12 -4 4 5 15 5 13 5 13 5 37 5 15 14 8 22 16 6 5 8 11 -2 9 1 6 23 10 7 3 12 14 7 5 7 7 5 5 37 8 20 39 14 9 6 7 7 14 11 9 6 9 7 7 11 21 12 5 26 0 8 31 16 8 4'

# --- ST sheet ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B4").Value = '59 69 62 68 64 55 73 69 61 64 68 62 64 65 62 64 64 65 60 63'
$stWs.Range("B5").Value = '14 22 18 19 35 15 27 13 22 23 26 20 22 27 27 17 16 20 18 20'
$stWs.Range("D3").Value = '58 50 36 58 51 60 52 59 45 47 43 59 47 41 13 47 52 50 44 57 55 52 43 49 49 62 55 52 53 51 51 46 51 47 50 41 43 48 33 41 58 52 34 49 49 39 31 43 47'
$stWs.Range("D4").Value = '0 0 0 18 0 -4 15 0 0 0 10 0 0 4 56 1 5 0 8 21 14 20 17 0 8 17 0 9 -2 27 0 0 10 0 16 0 6 3 0 8 0 0 0 10 14 0 0 17 -2'
$stWs.Range("D5").Value = '0 0 0 0 0 6 12 0 0 48 0 0 -5 0 0'
$stWs.Range("B2").Value = 87
$stWs.Range("D2").Value = 49
$stWs.Range("F2").Value = 549
$stWs.Range("G2").Value = 535
$stWs.Range("H2").Value = 3
$stWs.Range("J2").Value = 215
$stWs.Range("K2").Value = 202
$stWs.Range("N2").Value = 86
$stWs.Range("O2").Value = 65
$stWs.Range("B3").Value = 67

# --- OFF sheet ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value = 5
$offWs.Range("C2").Value = 218
$offWs.Range("D2").Value = 13
$offWs.Range("F2").Value = 43
$offWs.Range("G2").Value = 63
$offWs.Range("I2").Value = 14
$offWs.Range("J2").Value = 35
$offWs.Range("N2").Value = 15
$offWs.Range("O2").Value = 26
$offWs.Range("P2").Value = 15
$offWs.Range("C3").Value = 167
$offWs.Range("E3").Value = 36
$offWs.Range("F3").Value = 126
$offWs.Range("G3").Value = 49
$offWs.Range("I3").Value = 51
$offWs.Range("J3").Value = 54
$offWs.Range("L3").Value = 295
$offWs.Range("M3").Value = 212
$offWs.Range("Q3").Value = 608

# --- DEF sheet ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 190
$defWs.Range("E2").Value = 12
$defWs.Range("F2").Value = 61
$defWs.Range("G2").Value = 46
$defWs.Range("H2").Value = 5
$defWs.Range("J2").Value = 27
$defWs.Range("N2").Value = 29
$defWs.Range("O2").Value = 30
$defWs.Range("P2").Value = 14
$defWs.Range("B3").Value = 18
$defWs.Range("C3").Value = 167
$defWs.Range("E3").Value = 38
$defWs.Range("F3").Value = 92
$defWs.Range("G3").Value = 35
$defWs.Range("H3").Value = 31
$defWs.Range("I3").Value = 50
$defWs.Range("J3").Value = 56
$defWs.Range("L3").Value = 314
$defWs.Range("M3").Value = 203
$defWs.Range("Q3").Value = 559

# --- TURNS sheet ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("C3").Value = 12
$turnsWs.Range("E3").Value = 9

# --- PEN sheet ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 26
$penWs.Range("D2").Value = 7
$penWs.Range("B3").Value = 19
$penWs.Range("D4").Value = 9

